$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.497586578130722
$ws.Range("B1").Value = 0.8091548681259155
$ws.Range("C1").Value = 1.052488803863525
$ws.Range("D1").Value = 4.733420848846436
$ws.Range("E1").Value = 1.731515645980835
